# Update latest output (run 72)
$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet: rows 2 and 3 (A:F) ---
$wsSchedule.Range("A2").Value = 46040.29166666666
$wsSchedule.Range("B2").Value = 46040.79166666666
$wsSchedule.Range("C2").Value = 12
$wsSchedule.Range("D2").Value = 45.36
$wsSchedule.Range("E2").Value = -27.44120924999999
$wsSchedule.Range("F2").Value = -0.6049649305555554

$wsSchedule.Range("A3").Value = 46040.83333333334
$wsSchedule.Range("B3").Value = 46041
$wsSchedule.Range("C3").Value = 4
$wsSchedule.Range("D3").Value = 15.12
$wsSchedule.Range("E3").Value = 379.85008425
$wsSchedule.Range("F3").Value = 25.12236006944445

# --- Detailed sheet ---
$wsDetailed.Range("E2").Value = "OFF"
$wsDetailed.Range("E3").Value = "OFF"
$wsDetailed.Range("E4").Value = "OFF"

$wsDetailed.Range("B5").Value = 56.98
$wsDetailed.Range("E5").Value = "OFF"

$wsDetailed.Range("B6").Value = 56.97996
$wsDetailed.Range("E6").Value = "OFF"

$wsDetailed.Range("B7").Value = 36.2
$wsDetailed.Range("C7").Value = "historical"
$wsDetailed.Range("E7").Value = "OFF"

$wsDetailed.Range("B8").Value = 56.97996
$wsDetailed.Range("C8").Value = "historical"
$wsDetailed.Range("E8").Value = "OFF"

$wsDetailed.Range("B9").Value = 56.98
$wsDetailed.Range("E9").Value = "OFF"

$wsDetailed.Range("B15").Value = 50.98493
$wsDetailed.Range("B16").Value = 32.29746
$wsDetailed.Range("B17").Value = 11.61977
$wsDetailed.Range("B18").Value = 0
$wsDetailed.Range("B19").Value = 2.80927
$wsDetailed.Range("B20").Value = 0.51
$wsDetailed.Range("B22").Value = -4.81333
$wsDetailed.Range("B23").Value = 0.51
$wsDetailed.Range("B24").Value = 0.62873
$wsDetailed.Range("B25").Value = -4.81902
$wsDetailed.Range("B26").Value = 0
$wsDetailed.Range("B27").Value = -0.89155
$wsDetailed.Range("B28").Value = -5.69904
$wsDetailed.Range("B30").Value = -9.99
$wsDetailed.Range("B31").Value = -20.40712
$wsDetailed.Range("B32").Value = -12.43825
$wsDetailed.Range("B33").Value = -12.01
$wsDetailed.Range("B34").Value = -6.71741
$wsDetailed.Range("B35").Value = -6.51053
$wsDetailed.Range("B37").Value = 0.66214
$wsDetailed.Range("B38").Value = 3.98497
$wsDetailed.Range("B39").Value = 17.67251
$wsDetailed.Range("B40").Value = 45.74509

$wsDetailed.Range("B42").Value = 64.8901
$wsDetailed.Range("E42").Value = "ON"
$wsDetailed.Range("E43").Value = "ON"
$wsDetailed.Range("E44").Value = "ON"

$wsDetailed.Range("B45").Value = 46.2146
$wsDetailed.Range("E45").Value = "ON"

$wsDetailed.Range("B46").Value = 56.34493
$wsDetailed.Range("E46").Value = "ON"

$wsDetailed.Range("E47").Value = "ON"
$wsDetailed.Range("E48").Value = "ON"
$wsDetailed.Range("E49").Value = "ON"
